# daily auto push: 2026-01-09 22:35 UTC
# Insert a new data row (2026/01/10, 土, 6, 201) before the existing row 603,
# shifting row 603 and everything below it down by one row
# (old row 603..644 become new row 604..645).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 603 (shifts rows 603:644 down to 604:645)
$ws.Rows.Item(603).Insert()

# Populate the newly inserted row 603 with the new reading.
# Dates in this sheet are stored as plain text (e.g. "2026/01/10"), not as
# real Excel dates, so force text by prefixing with an apostrophe.
$ws.Cells.Item(603, 1).Value = "'2026/01/10"
$ws.Cells.Item(603, 2).Value = "土"
$ws.Cells.Item(603, 3).Value = 6
$ws.Cells.Item(603, 4).Value = 201
